$wb = $excel.ActiveWorkbook

# This script applies an automated data refresh to the per-job "Profits" sheets
# (ALC, ARM, BSM, CUL, GSM, LTW, WVR) -- updating market-price-derived columns
# H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
# K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ
# for specific leve rows, as produced by the scheduled market-data runner.

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 331.81818
$ws.Range("I39").Value = 80.125
$ws.Range("J39").Value = 1003
$ws.Range("K39").Value = 240.375
$ws.Range("L39").Value = 3009
$ws.Range("M39").Value = 55.625
$ws.Range("N39").Value = -3601
# Row 40
$ws.Range("H40").Value = 857.7273
$ws.Range("I40").Value = 711.3333
$ws.Range("J40").Value = 1171.4286
$ws.Range("K40").Value = 711.3333
$ws.Range("L40").Value = 1171.4286
$ws.Range("M40").Value = -536.3333
$ws.Range("N40").Value = -1521.4286
# Row 41
$ws.Range("H41").Value = 724.0714
$ws.Range("J41").Value = 886.9091
$ws.Range("L41").Value = 886.9091
$ws.Range("N41").Value = -1766.9091
# Row 53
$ws.Range("H53").Value = 1988.8
$ws.Range("I53").Value = 10
$ws.Range("J53").Value = 2483.5
$ws.Range("K53").Value = 10
$ws.Range("L53").Value = 2483.5
$ws.Range("M53").Value = 627
$ws.Range("N53").Value = -3757.5
# Row 64
$ws.Range("H64").Value = 3461.25
$ws.Range("I64").Value = 2768.5715
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 2768.5715
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2520.5715
$ws.Range("N64").Value = -4496
# Row 67
$ws.Range("H67").Value = 3461.25
$ws.Range("I67").Value = 2768.5715
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 2768.5715
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -1910.5715
$ws.Range("N67").Value = -5716
# Row 76
$ws.Range("H76").Value = 3270848
$ws.Range("I76").Value = 3067.375
$ws.Range("K76").Value = 3067.375
$ws.Range("M76").Value = -2752.375
# Row 79
$ws.Range("H79").Value = 3270848
$ws.Range("I79").Value = 3067.375
$ws.Range("K79").Value = 3067.375
$ws.Range("M79").Value = -1975.375
# Row 86
$ws.Range("H86").Value = 8853
$ws.Range("I86").Value = 1690
$ws.Range("K86").Value = 1690
$ws.Range("M86").Value = -567
# Row 89
$ws.Range("H89").Value = 8853
$ws.Range("I89").Value = 1690
$ws.Range("K89").Value = 8450
$ws.Range("M89").Value = -2834
# Row 112
$ws.Range("H112").Value = 4167733
$ws.Range("J112").Value = 1106.7858
$ws.Range("L112").Value = 3320.3574
$ws.Range("N112").Value = -5536.357400000001
# Row 116
$ws.Range("H116").Value = 15684815
$ws.Range("I116").Value = 35283572
$ws.Range("K116").Value = 35283572
$ws.Range("M116").Value = -35280130
# Row 132
$ws.Range("H132").Value = 38465132
$ws.Range("I132").Value = 45458524
$ws.Range("K132").Value = 136375572
$ws.Range("M132").Value = -136373042
# Row 138
$ws.Range("H138").Value = 3350.1333
$ws.Range("I138").Value = 2895.7693
$ws.Range("J138").Value = 3534.7188
$ws.Range("K138").Value = 8687.3079
$ws.Range("L138").Value = 10604.1564
$ws.Range("M138").Value = -3547.3079
$ws.Range("N138").Value = -20884.1564

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9317.074000000001
$ws.Range("I32").Value = 6623.603
$ws.Range("J32").Value = 23406
$ws.Range("K32").Value = 6623.603
$ws.Range("L32").Value = 23406
$ws.Range("M32").Value = -6336.603
$ws.Range("N32").Value = -23980
# Row 45
$ws.Range("H45").Value = 3161.52
$ws.Range("I45").Value = 2700.5264
$ws.Range("J45").Value = 4621.3335
$ws.Range("K45").Value = 2700.5264
$ws.Range("L45").Value = 4621.3335
$ws.Range("M45").Value = -2323.5264
$ws.Range("N45").Value = -5375.3335

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""
# Row 29
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = -10578
# Row 86
$ws.Range("H86").Value = 2688.6875
$ws.Range("I86").Value = 2393.7144
$ws.Range("J86").Value = 4753.5
$ws.Range("K86").Value = 2393.7144
$ws.Range("L86").Value = 4753.5
$ws.Range("M86").Value = -1270.7144
$ws.Range("N86").Value = -6999.5
# Row 89
$ws.Range("H89").Value = 2688.6875
$ws.Range("I89").Value = 2393.7144
$ws.Range("J89").Value = 4753.5
$ws.Range("K89").Value = 11968.572
$ws.Range("L89").Value = 23767.5
$ws.Range("M89").Value = -6352.572
$ws.Range("N89").Value = -34999.5

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 52.77778
$ws.Range("J33").Value = 69.166664
$ws.Range("L33").Value = 414.999984
$ws.Range("N33").Value = -980.999984
# Row 48
$ws.Range("H48").Value = 100
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
# Row 80
$ws.Range("H80").Value = 19333.5
$ws.Range("J80").Value = 28000.25
$ws.Range("L80").Value = 84000.75
$ws.Range("N80").Value = -85872.75
# Row 83
$ws.Range("H83").Value = 19333.5
$ws.Range("J83").Value = 28000.25
$ws.Range("L83").Value = 252002.25
$ws.Range("N83").Value = -261362.25
# Row 109
$ws.Range("H109").Value = 3019.8823
$ws.Range("I109").Value = 835.6
$ws.Range("J109").Value = 3552.634
$ws.Range("K109").Value = 2506.8
$ws.Range("L109").Value = 10657.902
$ws.Range("M109").Value = -1466.8
$ws.Range("N109").Value = -12737.902
# Row 117
$ws.Range("H117").Value = 1439.7778
$ws.Range("I117").Value = 994
$ws.Range("K117").Value = 2982
$ws.Range("M117").Value = 460
# Row 131
$ws.Range("H131").Value = 661.73
$ws.Range("J131").Value = 693.5281
$ws.Range("L131").Value = 2080.5843
$ws.Range("N131").Value = -12160.5843
# Row 137
$ws.Range("H137").Value = 15878898
$ws.Range("J137").Value = 18525130
$ws.Range("L137").Value = 55575390
$ws.Range("N137").Value = -55585590

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1710
# Row 70
$ws.Range("H70").Value = 3293844.8
$ws.Range("I70").Value = 4533.222
$ws.Range("K70").Value = 4533.222
$ws.Range("M70").Value = -4263.222
# Row 73
$ws.Range("H73").Value = 3293844.8
$ws.Range("I73").Value = 4533.222
$ws.Range("K73").Value = 4533.222
$ws.Range("M73").Value = -3597.222
# Row 80
$ws.Range("H80").Value = 3767.0833
$ws.Range("I80").Value = 3390.4
$ws.Range("J80").Value = 4036.1428
$ws.Range("K80").Value = 3390.4
$ws.Range("L80").Value = 4036.1428
$ws.Range("M80").Value = -2392.4
$ws.Range("N80").Value = -6032.1428
# Row 83
$ws.Range("H83").Value = 3767.0833
$ws.Range("I83").Value = 3390.4
$ws.Range("J83").Value = 4036.1428
$ws.Range("K83").Value = 16952
$ws.Range("L83").Value = 20180.714
$ws.Range("M83").Value = -11960
$ws.Range("N83").Value = -30164.714

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2624.75
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751
# Row 71
$ws.Range("H71").Value = 2624.75
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756
# Row 122
$ws.Range("H122").Value = 2454680.5
$ws.Range("I122").Value = 3270740.8
$ws.Range("K122").Value = 9812222.399999999
$ws.Range("M122").Value = -9809772.399999999

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 1026
$ws.Range("I32").Value = 1026
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1026
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -709
$ws.Range("N32").Value = ""
# Row 136
$ws.Range("H136").Value = 32261540
$ws.Range("I136").Value = 47620868
$ws.Range("K136").Value = 142862604
$ws.Range("M136").Value = -142860054
